$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Remove the "Test Report" sheet (summary/report sheet no longer needed) ---
$wsReport = $wb.Sheets.Item("Test Report")
$wsReport.Delete()

# --- Rename the remaining sheet to "Sprint 3" ---
$ws = $wb.Sheets.Item("Patient-Clinical Data")
$ws.Name = "Sprint 3"

# --- Update test-run dates that were re-verified later (TC1, TC2, TC5 -> Mar 13 2018; TC6 -> Mar 11 2018) ---
$ws.Range("H8").Value = 43172
$ws.Range("H9").Value = 43172
$ws.Range("H12").Value = 43172
$ws.Range("H13").Value = 43170

# --- Add the new TC7 test case row (row 14) ---
$ws.Range("B14").Value = "Restaurants nearby:`nShow nearby restaurants close to the current user location (i.e. Centennial College)"
$ws.Range("C14").Value = "Open the application and click / tap onto a restaurant nearby the user location (i.e. Centennial College) in the map."
$ws.Range("D14").Value = "The Restaurant details (i.e. restaurant name appears)"
$ws.Range("G14").Value = "Fernando"
$ws.Range("H14").Value = 43172
$ws.Range("I14").Value = "Pass"
$ws.Rows.Item(14).RowHeight = 72

# --- Move selection to the newly added row ---
$ws.Range("H14").Select()
